$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin -> Bitcoin
$ws.Cells.Item(2, 4).Value = "56.705.57"
$ws.Cells.Item(2, 5).Value = "  +1.52%  "

# Row 3: Ethereum -> Ethereum
$ws.Cells.Item(3, 4).Value = "2.485.79"
$ws.Cells.Item(3, 5).Value = "  -1.36%  "

# Row 4: TetherUSD -> TetherUSD
$ws.Cells.Item(4, 5).Value = "  +0.10%  "

# Row 5: BNB -> BNB
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "488.69"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.90%  "

# Row 6: Solana -> Solana
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "148.77"
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +6.99%  "

# Row 7: USDC -> USDC
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.997"
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.19%  "

# Row 8: XRP -> XRP
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.511"
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -0.01%  "

# Row 9: LidoStakedEther -> LidoStakedEther
$ws.Cells.Item(9, 4).Value = "2.494.09"
$ws.Cells.Item(9, 5).Value = "  -1.03%  "

# Row 10: Toncoin -> Toncoin
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.77"
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +5.51%  "

# Row 11: Dogecoin -> Dogecoin
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0983"
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -1.19%  "

# Row 12: Cardano -> Cardano
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.335"
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +2.26%  "

# Row 13: TRON -> TRON
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.124"
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +1.26%  "

# Row 14: WrappedliquidstakedEther2.0 -> WrappedliquidstakedEther2.0
$ws.Cells.Item(14, 4).Value = "2.918.51"
$ws.Cells.Item(14, 5).Value = "  -1.19%  "

# Row 15: WrappedBTC -> WrappedBTC
$ws.Cells.Item(15, 4).Value = "56.582.39"
$ws.Cells.Item(15, 5).Value = "  +1.41%  "

# Row 16: Avalanche -> Avalanche
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "21.15"
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +2.32%  "

# Row 17: ShibaInu -> ShibaInu
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0000137"
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -0.98%  "

# Row 18: WrappedEther -> WrappedEther
$ws.Cells.Item(18, 4).Value = "2.478.97"
$ws.Cells.Item(18, 5).Value = "  -1.54%  "

# Row 19: Polkadot -> Polkadot
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.57"
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +4.64%  "

# Row 20: Chainlink -> Chainlink
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "10.26"
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +2.14%  "

# Row 21: BitcoinCash -> BitcoinCash
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "319.67"
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -1.34%  "

# Row 22: Dai -> Dai
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.03%  "

# Row 23: Uniswap -> Uniswap
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.92"
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +3.55%  "

# Row 24: Litecoin -> Litecoin
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "58.37"
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.69%  "

# Row 25: Polygon -> Polygon
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.415"
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +2.96%  "

# Row 26: Binance-PegBSC-USD -> Binance-PegBSC-USD
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.997"
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -1.29%  "

# Row 27: Kaspa -> Kaspa
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.163"
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.30%  "

# Row 28: WrappedeETH -> WrappedeETH
$ws.Cells.Item(28, 4).Value = "2.567.37"
$ws.Cells.Item(28, 5).Value = "  -1.66%  "

# Row 29: InternetComputer(DFINITY) -> InternetComputer(DFINITY)
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.73"
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +4.46%  "

# Row 30: PEPE -> PEPE
$ws.Cells.Item(30, 4).Value = "0.0₃0790"
$ws.Cells.Item(30, 5).Value = "  +3.02%  "

# Row 31: USDe -> USDe
$ws.Cells.Item(31, 5).Value = "  +0.00%  "

# Row 32: Monero -> Monero
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "149.08"
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -0.97%  "

# Row 33: EthereumClassic -> EthereumClassic
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "18.34"
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +0.79%  "

# Row 34: PancakeSwap -> PancakeSwap
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.51"
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +1.57%  "

# Row 35: Aptos -> Aptos
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.21"
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -0.29%  "

# Row 36: ImmutableX -> ImmutableX
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.17"
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +6.94%  "

# Row 37: NEARProtocol -> NEARProtocol
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.77"
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +1.54%  "

# Row 38: Fetch.AI -> Fetch.AI
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.874"
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +4.73%  "

# Row 39: OKB -> OKB
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "34.19"
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -1.30%  "

# Row 40: Filecoin -> Stacks
$ws.Cells.Item(40, 2).Value = "Stacks"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.38"
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +3.97%  "

# Row 41: FirstDigitalUSD -> Filecoin
$ws.Cells.Item(41, 2).Value = "Filecoin"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.52"
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +4.02%  "

# Row 42: Hedera -> Hedera
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0557"
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +1.88%  "

# Row 43: Mantle -> FirstDigitalUSD
$ws.Cells.Item(43, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.995"
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.43%  "

# Row 44: Stacks -> Mantle
$ws.Cells.Item(44, 2).Value = "Mantle"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.610"
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.82%  "

# Row 45: Stellar -> Bittensor
$ws.Cells.Item(45, 2).Value = "Bittensor"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "265.31"
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +7.77%  "

# Row 46: Bittensor -> RenderToken
$ws.Cells.Item(46, 2).Value = "RenderToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.76"
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +8.19%  "

# Row 47: WhiteBITCoin -> VeChain
$ws.Cells.Item(47, 2).Value = "VeChain"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0232"
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +4.17%  "

# Row 48: VeChain -> Stellar
$ws.Cells.Item(48, 2).Value = "Stellar"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0928"
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +3.04%  "

# Row 49: RenderToken -> WhiteBITCoin
$ws.Cells.Item(49, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "10.22"
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +0.82%  "

# Row 50: EnergySwap -> EnergySwap
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "17.74"
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.75%  "

# Row 51: Maker -> Maker
$ws.Cells.Item(51, 4).Value = "1.890.99"
$ws.Cells.Item(51, 5).Value = "  -4.50%  "

